$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.42638983778676
$ws.Range("C2").Value = 10.29331459755666
$ws.Range("D2").Value = 7.30167665320554
$ws.Range("E2").Value = 16.45433944982152
$ws.Range("F2").Value = 43.42741042073344
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("B3").Value = 15.82553249713264
$ws.Range("C3").Value = 9.688831710892437
$ws.Range("D3").Value = 7.136393549518607
$ws.Range("E3").Value = 15.51418462611142
$ws.Range("F3").Value = 41.92871716270902
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("B4").Value = 15.45323296866087
$ws.Range("C4").Value = 9.302815878495984
$ws.Range("D4").Value = 7.033694688316622
$ws.Range("E4").Value = 14.91370740114557
$ws.Range("F4").Value = 40.99073186571778
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("B5").Value = 15.30097127547397
$ws.Range("C5").Value = 9.141972636156947
$ws.Range("D5").Value = 6.991576212751836
$ws.Range("E5").Value = 14.66344036576686
$ws.Range("F5").Value = 40.60455364992115
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("B6").Value = 15.27566457246559
$ws.Range("C6").Value = 9.115057593640415
$ws.Range("D6").Value = 6.984567413363458
$ws.Range("E6").Value = 14.62155641059596
$ws.Range("F6").Value = 40.54020732386432
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("B7").Value = 15.4511813125667
$ws.Range("C7").Value = 9.300660726981791
$ws.Range("D7").Value = 7.033127697937479
$ws.Range("E7").Value = 14.91035435876094
$ws.Range("F7").Value = 40.98553895853991
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("B8").Value = 16.22008705086596
$ws.Range("C8").Value = 10.08808041295122
$ws.Range("D8").Value = 7.244961183216043
$ws.Range("E8").Value = 16.13513958936087
$ws.Range("F8").Value = 42.91465571275551
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("B9").Value = 17.68917423310451
$ws.Range("C9").Value = 11.65565728110506
$ws.Range("D9").Value = 7.649156145900511
$ws.Range("E9").Value = 18.42300014397892
$ws.Range("F9").Value = 46.53555840643371
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("B10").Value = 18.73062856694719
$ws.Range("C10").Value = 12.72416089805899
$ws.Range("D10").Value = 7.937297104013416
$ws.Range("E10").Value = 20.09240612184273
$ws.Range("F10").Value = 49.07208829654715
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("B11").Value = 19.19367544082326
$ws.Range("C11").Value = 13.18385019090909
$ws.Range("D11").Value = 8.066081172719578
$ws.Range("E11").Value = 20.81095225300755
$ws.Range("F11").Value = 50.19475926645455
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("B12").Value = 19.36730150894325
$ws.Range("C12").Value = 13.3541507046749
$ws.Range("D12").Value = 8.114489260347284
$ws.Range("E12").Value = 21.07721539202315
$ws.Range("F12").Value = 50.61508653765809
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("B13").Value = 19.32998712384931
$ws.Range("C13").Value = 13.31764096272882
$ws.Range("D13").Value = 8.104080161986687
$ws.Range("E13").Value = 21.02012967387631
$ws.Range("F13").Value = 50.52477981594377
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("B14").Value = 19.20799517259175
$ws.Range("C14").Value = 13.19793643458213
$ws.Range("D14").Value = 8.07007110694499
$ws.Range("E14").Value = 20.83297457037248
$ws.Range("F14").Value = 50.22943779388444
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("B15").Value = 19.13304269629683
$ws.Range("C15").Value = 13.12412308575819
$ws.Range("D15").Value = 8.049191890136083
$ws.Range("E15").Value = 20.7175780972205
$ws.Range("F15").Value = 50.04789795263731
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("B16").Value = 18.70013593250552
$ws.Range("C16").Value = 12.69358994218468
$ws.Range("D16").Value = 7.928832017558936
$ws.Range("E16").Value = 20.04462861492908
$ws.Range("F16").Value = 48.99806319788954
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("B17").Value = 18.43168260467048
$ws.Range("C17").Value = 12.42272897323808
$ws.Range("D17").Value = 7.854385926304401
$ws.Range("E17").Value = 19.62135747802872
$ws.Range("F17").Value = 48.34579288210324
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("B18").Value = 18.27627591557962
$ws.Range("C18").Value = 12.26445279453743
$ws.Range("D18").Value = 7.811352361224981
$ws.Range("E18").Value = 19.37405318474894
$ws.Range("F18").Value = 47.96770302184997
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("B19").Value = 18.22349204673524
$ws.Range("C19").Value = 12.21043598561017
$ws.Range("D19").Value = 7.796746104326484
$ws.Range("E19").Value = 19.28965765983189
$ws.Range("F19").Value = 47.83919683658041
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("B20").Value = 18.46036472424897
$ws.Range("C20").Value = 12.45181954634602
$ws.Range("D20").Value = 7.862333226540804
$ws.Range("E20").Value = 19.6668135541616
$ws.Range("F20").Value = 48.41553291880169
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("B21").Value = 19.24387515290404
$ws.Range("C21").Value = 13.23319881649802
$ws.Range("D21").Value = 8.080070379528696
$ws.Range("E21").Value = 20.88810456907759
$ws.Range("F21").Value = 50.3163196081298
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("B22").Value = 19.74585311027923
$ws.Range("C22").Value = 13.72188835381062
$ws.Range("D22").Value = 8.220265545473767
$ws.Range("E22").Value = 21.65229987435612
$ws.Range("F22").Value = 51.53047303727401
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("B23").Value = 19.47891475148388
$ws.Range("C23").Value = 13.46307063694064
$ws.Range("D23").Value = 8.145643027621672
$ws.Range("E23").Value = 21.24753019877283
$ws.Range("F23").Value = 50.8851241141432
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("B24").Value = 18.4474008547955
$ws.Range("C24").Value = 12.43867565670082
$ws.Range("D24").Value = 7.858740977155679
$ws.Range("E24").Value = 19.64627520016255
$ws.Range("F24").Value = 48.38401305149029
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("B25").Value = 17.29746576781898
$ws.Range("C25").Value = 11.23866050126309
$ws.Range("D25").Value = 7.541222888398656
$ws.Range("E25").Value = 17.77169398848133
$ws.Range("F25").Value = 45.57610727184977
$ws.Range("H25").Value = 7.344005520526261
